$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "1.007") but must remain
# plain text, matching the workbook's original inline-string representation.
# Force a Text number format first so Excel does not auto-convert the string to a
# numeric value, then clear the format afterwards so the cell keeps the same
# (unstyled) appearance it had before the edit.
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D42",
    "D43",
    "D46",
    "D48",
    "D50",
    "D51"
)

foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated coin prices (column D) and hourly volume deltas (column E)
$ws.Range("D2").Value = "26.247.81"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.684.17"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "217.68"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "0.5248"
$ws.Range("E6").Value = "  +2.37%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.2703"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "0.06410"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").Value = "0.07484"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "1.708.57"
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "4.546"
$ws.Range("D14").Value = "0.5804"
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("D15").Value = "0.000008466"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "64.25"
$ws.Range("E16").Value = "  -0.99%  "
$ws.Range("D17").Value = "26.310.92"
$ws.Range("D18").Value = "4.919"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "10.86"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "189.24"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "6.196"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "144.31"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Value = "7.693"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "0.1238"
$ws.Range("E26").Value = "  +4.63%  "
$ws.Range("D27").Value = "15.77"
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").Value = "0.06652"
$ws.Range("E28").Value = "  +12.26%  "
$ws.Range("D29").Value = "1.346"
$ws.Range("E29").Value = "  +5.23%  "
$ws.Range("D30").Value = "1.327"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "3.569"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").Value = "3.563"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "1.664"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("D35").Value = "0.6196"
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").Value = "2.397"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").Value = "2.706"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").Value = "6.388"
$ws.Range("E38").Value = "  +5.53%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").Value = "1.104.67"
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("D42").Value = "1.015"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "100.53"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "1.831.90"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "56.69"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "8.152"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  +1.36%  "
$ws.Range("D50").Value = "0.4302"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "6.022"
$ws.Range("E51").Value = "  +2.35%  "

# Restore the default (unstyled) formatting on the cells we forced to Text above
foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
